# "fixed typo in data dictionary"
#
# The "Analysis" column (G) uses a lowercase "x" marker everywhere except a
# handful of rows that had accidentally been marked with an uppercase "X".
# This script normalizes those rows to the lowercase marker, fills in the
# (until-now blank) Analysis column for every other analysis-variable row,
# removes the erroneous mark on row 46 (which is not an analysis variable),
# and corrects three Description cells (rows 79/82/85) that had been
# copy-pasted from the "Driving time" variables instead of the "Walking
# time" ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows where column G ("Analysis") should read lowercase "x".
# (Rows 42 and 43 already had an uppercase "X" that gets corrected here;
# the rest are rows that had no mark at all.)
$xRows = @(7,11,12,13,14,15,16,17,18,19,33,34,35,36,37,38,42,43,44,45,48,50,55,56,63,68,117,119,120,121,122)

foreach ($r in $xRows) {
    $ws.Cells.Item($r, 7).Value = "x"
}

# Row 46 incorrectly had an "X" Analysis mark; remove it entirely.
$ws.Cells.Item(46, 7).Clear()

# Fix copy/paste typo in the Description column (I) for the walking-time
# variables, which had been pointing at the driving-time descriptions.
$ws.Cells.Item(79, 9).Value = "Walking time (minutes) to nearest buprenorphine provider"
$ws.Cells.Item(82, 9).Value = "Walking time (minutes) to nearest methadone provider"
$ws.Cells.Item(85, 9).Value = "Walking time (minutes) to nearest naltrexone provider"

# Restore the cursor/selection to where the editor left it.
$ws.Range("G18").Select()
